# chore: update Sheets via scheduled runner
# Refresh cached Universalis price snapshots (currentAveragePrice*, LevePrice*,
# LeveProfit*) for the affected leves across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 137
$ws.Cells.Item(137, 8).Value = 1454.1111
$ws.Cells.Item(137, 10).Value = 2229.7144
$ws.Cells.Item(137, 12).Value = 6689.1432
$ws.Cells.Item(137, 14).Value = -11789.1432

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).Value = ""

# Row 32
$ws.Cells.Item(32, 8).Value = 3464.2307
$ws.Cells.Item(32, 9).Value = 2774.855
$ws.Cells.Item(32, 11).Value = 2774.855
$ws.Cells.Item(32, 13).Value = -2487.855

# Row 47
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(47, 12).Value = 0
$ws.Cells.Item(47, 14).Value = ""

# Row 101
$ws.Cells.Item(101, 8).Value = 43386.2
$ws.Cells.Item(101, 10).Value = 43386.2
$ws.Cells.Item(101, 12).Value = 43386.2
$ws.Cells.Item(101, 14).Value = -49876.2

# Row 132
$ws.Cells.Item(132, 8).Value = 1318.6171
$ws.Cells.Item(132, 9).Value = 948.1
$ws.Cells.Item(132, 11).Value = 2844.3
$ws.Cells.Item(132, 13).Value = -314.3000000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Cells.Item(22, 8).Value = 356
$ws.Cells.Item(22, 9).Value = 300
$ws.Cells.Item(22, 10).Value = 370
$ws.Cells.Item(22, 11).Value = 300
$ws.Cells.Item(22, 12).Value = 370
$ws.Cells.Item(22, 13).Value = -127
$ws.Cells.Item(22, 14).Value = -716

# Row 86
$ws.Cells.Item(86, 8).Value = 121304.06
$ws.Cells.Item(86, 9).Value = 4324.5454
$ws.Cells.Item(86, 10).Value = 335766.5
$ws.Cells.Item(86, 11).Value = 4324.5454
$ws.Cells.Item(86, 12).Value = 335766.5
$ws.Cells.Item(86, 13).Value = -3201.5454
$ws.Cells.Item(86, 14).Value = -338012.5

# Row 89
$ws.Cells.Item(89, 8).Value = 121304.06
$ws.Cells.Item(89, 9).Value = 4324.5454
$ws.Cells.Item(89, 10).Value = 335766.5
$ws.Cells.Item(89, 11).Value = 21622.727
$ws.Cells.Item(89, 12).Value = 1678832.5
$ws.Cells.Item(89, 13).Value = -16006.727
$ws.Cells.Item(89, 14).Value = -1690064.5

# Row 99
$ws.Cells.Item(99, 8).Value = 2846.4
$ws.Cells.Item(99, 9).Value = 2861.5386
$ws.Cells.Item(99, 10).Value = 2748
$ws.Cells.Item(99, 11).Value = 2861.5386
$ws.Cells.Item(99, 12).Value = 2748
$ws.Cells.Item(99, 13).Value = -1363.5386
$ws.Cells.Item(99, 14).Value = -5744

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 2382572.5
$ws.Cells.Item(31, 10).Value = 1771.238
$ws.Cells.Item(31, 12).Value = 1771.238
$ws.Cells.Item(31, 14).Value = -2361.238

# Row 34
$ws.Cells.Item(34, 8).Value = 2382572.5
$ws.Cells.Item(34, 10).Value = 1771.238
$ws.Cells.Item(34, 12).Value = 1771.238
$ws.Cells.Item(34, 14).Value = -2175.238

# Row 53
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 14).Value = ""

# Row 58
$ws.Cells.Item(58, 8).Value = 1146438.6
$ws.Cells.Item(58, 9).Value = 1891950.4
$ws.Cells.Item(58, 11).Value = 1891950.4
$ws.Cells.Item(58, 13).Value = -1891747.4

# Row 134
$ws.Cells.Item(134, 8).Value = 1585.8605
$ws.Cells.Item(134, 9).Value = 1334.3948
$ws.Cells.Item(134, 10).Value = 3497
$ws.Cells.Item(134, 11).Value = 4003.1844
$ws.Cells.Item(134, 12).Value = 10491
$ws.Cells.Item(134, 13).Value = -1468.1844
$ws.Cells.Item(134, 14).Value = -15561

# Row 136
$ws.Cells.Item(136, 8).Value = 1146438.6
$ws.Cells.Item(136, 9).Value = 1891950.4
$ws.Cells.Item(136, 11).Value = 5675851.199999999
$ws.Cells.Item(136, 13).Value = -5673301.199999999

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 8680.674000000001
$ws.Cells.Item(131, 9).Value = 360.91666
$ws.Cells.Item(131, 10).Value = 9977.26
$ws.Cells.Item(131, 11).Value = 1082.74998
$ws.Cells.Item(131, 12).Value = 29931.78
$ws.Cells.Item(131, 13).Value = 3957.25002
$ws.Cells.Item(131, 14).Value = -40011.78

# Row 132
$ws.Cells.Item(132, 8).Value = 1171.375
$ws.Cells.Item(132, 10).Value = 1443.25
$ws.Cells.Item(132, 12).Value = 12989.25
$ws.Cells.Item(132, 14).Value = -18049.25

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 48.083332
$ws.Cells.Item(2, 9).Value = 11.076923
$ws.Cells.Item(2, 11).Value = 11.076923
$ws.Cells.Item(2, 13).Value = 101.923077

# Row 7
$ws.Cells.Item(7, 8).Value = 4091058.8
$ws.Cells.Item(7, 9).Value = 5375000
$ws.Cells.Item(7, 10).Value = 1009599.6
$ws.Cells.Item(7, 11).Value = 5375000
$ws.Cells.Item(7, 12).Value = 1009599.6
$ws.Cells.Item(7, 13).Value = -5374888
$ws.Cells.Item(7, 14).Value = -1009823.6

# Row 8
$ws.Cells.Item(8, 8).Value = 4091058.8
$ws.Cells.Item(8, 9).Value = 5375000
$ws.Cells.Item(8, 10).Value = 1009599.6
$ws.Cells.Item(8, 11).Value = 5375000
$ws.Cells.Item(8, 12).Value = 1009599.6
$ws.Cells.Item(8, 13).Value = -5374861
$ws.Cells.Item(8, 14).Value = -1009877.6

# Row 11
$ws.Cells.Item(11, 8).Value = 4637170.5
$ws.Cells.Item(11, 9).Value = 5441244
$ws.Cells.Item(11, 11).Value = 5441244
$ws.Cells.Item(11, 13).Value = -5441105

# Row 12
$ws.Cells.Item(12, 8).Value = 6142857
$ws.Cells.Item(12, 10).Value = 7000000
$ws.Cells.Item(12, 12).Value = 7000000
$ws.Cells.Item(12, 14).Value = -7000280

# Row 18
$ws.Cells.Item(18, 8).Value = 10000000
$ws.Cells.Item(18, 9).Value = 10000000
$ws.Cells.Item(18, 11).Value = 10000000
$ws.Cells.Item(18, 13).Value = -9999707

# Row 97
$ws.Cells.Item(97, 8).Value = 1183.0769
$ws.Cells.Item(97, 9).Value = 1031.6666
$ws.Cells.Item(97, 11).Value = 1031.6666
$ws.Cells.Item(97, 13).Value = -535.6666

# Row 123
$ws.Cells.Item(123, 8).Value = 15512.75
$ws.Cells.Item(123, 10).Value = 15512.75
$ws.Cells.Item(123, 12).Value = 15512.75
$ws.Cells.Item(123, 14).Value = -20412.75

# Row 132
$ws.Cells.Item(132, 8).Value = 963890.6
$ws.Cells.Item(132, 9).Value = 1426133.9
$ws.Cells.Item(132, 10).Value = 3847.1538
$ws.Cells.Item(132, 11).Value = 4278401.699999999
$ws.Cells.Item(132, 12).Value = 11541.4614
$ws.Cells.Item(132, 13).Value = -4275871.699999999
$ws.Cells.Item(132, 14).Value = -16601.4614

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 1917.8182
$ws.Cells.Item(46, 9).Value = 1025
$ws.Cells.Item(46, 11).Value = 1025
$ws.Cells.Item(46, 13).Value = -837

# Row 100
$ws.Cells.Item(100, 8).Value = 2425
$ws.Cells.Item(100, 9).Value = 2000
$ws.Cells.Item(100, 11).Value = 2000
$ws.Cells.Item(100, 13).Value = -1459

# Row 136
$ws.Cells.Item(136, 8).Value = 1652.8
$ws.Cells.Item(136, 9).Value = 1001.94
$ws.Cells.Item(136, 10).Value = 4907.1
$ws.Cells.Item(136, 11).Value = 3005.82
$ws.Cells.Item(136, 12).Value = 14721.3
$ws.Cells.Item(136, 13).Value = -455.8200000000002
$ws.Cells.Item(136, 14).Value = -19821.3

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 1173.0123
$ws.Cells.Item(132, 9).Value = 759.8280999999999
$ws.Cells.Item(132, 10).Value = 2728.5293
$ws.Cells.Item(132, 11).Value = 2279.4843
$ws.Cells.Item(132, 12).Value = 8185.5879
$ws.Cells.Item(132, 13).Value = 250.5156999999999
$ws.Cells.Item(132, 14).Value = -13245.5879

# Row 133
$ws.Cells.Item(133, 8).Value = 60000
$ws.Cells.Item(133, 10).Value = 60000
$ws.Cells.Item(133, 12).Value = 60000
$ws.Cells.Item(133, 14).Value = -70120

# Row 136
$ws.Cells.Item(136, 8).Value = 10289627
$ws.Cells.Item(136, 9).Value = 13229021
$ws.Cells.Item(136, 10).Value = 1747.5
$ws.Cells.Item(136, 11).Value = 39687063
$ws.Cells.Item(136, 12).Value = 5242.5
$ws.Cells.Item(136, 13).Value = -39684513
$ws.Cells.Item(136, 14).Value = -10342.5
